$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: shift "living_rooms_2" from E1 to B1, and push kitchens_1/bedrooms_1/bedrooms_2
# each one column to the right (B1->C1, C1->D1, D1->E1). A1 and F1 stay unchanged.
$ws.Range("B1").Value = "living_rooms_2"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "bedrooms_1"
$ws.Range("E1").Value = "bedrooms_2"

# Data rows: update the one-hot marker cells to match the new column order / new values
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 1

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 1

$ws.Range("A4").Value = 1
$ws.Range("E4").Value = 0

$ws.Range("A7").Value = 0
$ws.Range("B7").Value = 1
